$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.08737051486969
$ws.Range("B1").Value = 1.922170877456665
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.076799631118774
$ws.Range("E1").Value = 1.139383792877197
